$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 updates
$ws.Range("K2").Value = 1.92
$ws.Range("T2").Value = 1.5

# Row 4 updates
$ws.Range("G4").Value = 1.7
$ws.Range("H4").Value = 3.5
$ws.Range("I4").Value = 4.33
$ws.Range("J4").Value = 2.4
$ws.Range("L4").Value = 5.5
$ws.Range("M4").Value = 1.07
$ws.Range("N4").Value = 9
$ws.Range("O4").Value = 1.36
$ws.Range("P4").Value = 3
$ws.Range("AA4").Value = 2.1
$ws.Range("AB4").Value = 1.67
$ws.Range("AD4").Value = 7.5
$ws.Range("AF4").Value = 13
$ws.Range("AI4").Value = 8.5
$ws.Range("AJ4").Value = 7
$ws.Range("AK4").Value = 21
$ws.Range("AN4").Value = 10
$ws.Range("AO4").Value = 23
$ws.Range("AP4").Value = 15
$ws.Range("AQ4").Value = 51
$ws.Range("AR4").Value = 41
$ws.Range("AS4").Value = 51
